# Update the lab-related strings on the first worksheet.
# "VA.MHV.PHR.lab"   -> "VA.MHV.PHR.labTest"
# "VA MHV PHR lab"   -> "VA MHV PHR labTest"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "VA.MHV.PHR.labTest"
$ws.Range("B2").Value = "VA MHV PHR labTest"
